# Generate Report for Handoff
# Updates the localization-status report:
#  - "Ready for handoff" rows (8, 9, 11, 12, 13, 14) now have a Priority of "ht"
#    on the zh-cn and de-de sheets.
#  - The corresponding "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#    timestamps are refreshed to reflect the new handoff generation time.

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 11, 12, 13, 14)

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

foreach ($r in $rows) {
    # Mark these files as handed off with priority "ht" in both locale sheets.
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"

    # Refresh the handoff timestamps to the newly generated values.
    $zhcn.Range("H$r").Value = "2016-08-25 18:21:55"
    $dede.Range("H$r").Value = "2016-08-25 18:22:00"
    $overview.Range("G$r").Value = "2016-08-25 18:22:00"
}
